# Build RAD IMR 1.1.0
# Updates the "Metadata" sheet of the ValueSet workbook: bump Version,
# add Experimental flag, refresh build Date, split the Contact info into
# three rows, update Jurisdiction, and fix wording in Description.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Version: 1.0.0 -> 1.1.0
$ws.Range("B3").Value = "1.1.0"

# Experimental: previously blank, now the literal text "false".
# A bare Value = "false" auto-coerces to the Boolean FALSE (same as
# typing it straight into Excel), so build it as a formula result and
# paste it back in as a value to keep it stored as text.
$ws.Range("D1").Formula = '=""&"false"'
$ws.Range("D1").Copy()
$ws.Range("B7").PasteSpecial(-4163)
$ws.Range("D1").Value = ""

# Date: refreshed build timestamp
$ws.Range("B8").Value = "2024-06-20T08:51:57-05:00"

# Contact block now spans three rows (was a single "No display for
# ContactDetail" row) - rows 10-12 are all labelled "Contact" in column A.
$ws.Range("A11").Value = "Contact"
$ws.Range("A12").Value = "Contact"

$ws.Range("B10").Value = "null (https://www.ihe.net/ihe_domains/radiology/)"
$ws.Range("B11").Value = "null (radiology@ihe.net)"
$ws.Range("B12").Value = "IHE Radiology Technical Committee (radiology@ihe.net)"

# Jurisdiction: World -> Global (Whole world)
$ws.Range("B13").Value = "Global (Whole world)"

# Description: fix wording
$ws.Range("B14").Value = "Codes representing the applicable intent for an imaging ServiceRequest."
